$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed error in Total forest: update cached values in columns C (Other) and G (Drevet_Skov)
$ws.Range("C2").Value = -208.084300000001
$ws.Range("G2").Value = 629.8975

$ws.Range("C4").Value = 158.078100000001
$ws.Range("G4").Value = 486.1285

$ws.Range("C5").Value = 48.2409
$ws.Range("G5").Value = 85.6039

$ws.Range("C6").Value = 11.4981
$ws.Range("G6").Value = 135.8964

$ws.Range("C8").Value = 80.7752
$ws.Range("G8").Value = 4.1169

$ws.Range("C9").Value = 1.69289999999995
$ws.Range("G9").Value = 90.0854

$ws.Range("C10").Value = -56.5546
$ws.Range("G10").Value = 25.0017

$ws.Range("C11").Value = 3.5311
$ws.Range("G11").Value = 17.105
